$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-annotated dialog acts (DAMSLTag in column I, DialogAct in column J)
$updates = @(
    @{ Row = 5;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 28; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 40; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 47; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 48; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 50; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 55; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 61; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 63; Tag = "ba"; Act = "Appreciation" }
)

foreach ($u in $updates) {
    $ws.Range("I$($u.Row)").Value = $u.Tag
    $ws.Range("J$($u.Row)").Value = $u.Act
}
